# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as the new row 3 of the daily
# logic subset sheet; every existing data row (previously rows 3-48) shifts
# down by one (to rows 4-49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 3 - this pushes rows 3..48 down
# to 4..49 (and grows the used range to A1:R49), matching Excel's normal
# "insert sheet row" behaviour (including carrying the date-format style
# of the row below down onto the new blank row).
$ws.Range("A3").EntireRow.Insert()

# Populate the newly inserted row 3 with the new observation.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44812
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12500
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 500
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
